# Append new benchmark result rows (rows 84-95) to the "Results" sheet,
# extending the used range from A1:G83 to A1:G95.
#
# Columns: A=runs(n), B=instance(text), C=encoding(text), D=time(text,
# even when it looks numeric, e.g. "0.031"), E=status(text),
# F=vars(n), G=clauses(n)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Results")

$newData = @(
    @(1, "5-4-4",  "binary",   "0.031",   "sat",     4161,  34227),
    @(1, "3-2-2",  "binary",   "0.000",   "sat",     168,   434),
    @(1, "3-2-2",  "binary",   "0.000",   "sat",     168,   434),
    @(1, "3-2-2",  "binary",   "0.000",   "sat",     168,   434),
    @(1, "3-2-2",  "binary",   "0.000",   "sat",     169,   480),
    @(1, "3-2-2",  "binary",   "0.000",   "sat",     169,   480),
    @(1, "3-2-2",  "binary",   "0.000",   "sat",     169,   480),
    @(1, "3-2-2",  "binary",   "0.000",   "sat",     169,   480),
    @(1, "3-2-2",  "binary",   "0.000",   "sat",     169,   480),
    @(1, "6-5-5",  "binary",   "timeout", "timeout", 12002, 160072),
    @(1, "6-5-5",  "binomial", "33.672",  "sat",     4500,  167615),
    @(1, "6-5-5",  "binary",   "9.625",   "sat",     12002, 157640)
)

$startRow = 84

# Column D values are always plain text in this sheet (even values that
# look like numbers, e.g. "0.031" or "33.672"), so force the whole
# destination range in column D to text format before writing values,
# which prevents Excel from auto-converting numeric-looking strings.
$lastRow = $startRow + $newData.Count - 1
$ws.Range("D$startRow`:D$lastRow").NumberFormat = "@"

$r = $startRow
foreach ($row in $newData) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $r++
}
